# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" rows for worker 73435907 (JHON JAIRO
# CASTILLO SABACH) were stored most-recent-first (2002 .. 1908). This
# update re-sorts them chronologically ascending (1908 .. 2002), keeping
# each period's corresponding "Valor Mora" (F column) attached to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Periods (column E) and their matching mora values (column F) for rows 18-24,
# in the new chronological order.
$periods = @("1908", "1909", "1910", "1911", "1912", "2001", "2002")
$valores = @(33125, 33125, 33125, 33125, 33125, 33125, 18771)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 18 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}

$wb.Save()
